# fix(excel), #16: Ranges are not valid candidate to unblock a cycle.
#
# Adds four new data rows (8-11, plus a trailing blank row 12 that only
# carries formatting) to the "circular.xlsx" test workbook so that it also
# exercises array-formula / range candidates when looking for a cell to
# "break" a circular reference.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 8: A8 is TRUE, D8 is an array formula (range) feeding B8/C8 back in
# a circular loop - this is the "array/range" candidate the fix targets.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = $true
$ws.Range("B8").Formula = "=C8"
$ws.Range("C8").Formula = "=D8"
$ws.Range("D8:D9").FormulaArray = "=IF(A8:A9,B8:B9,1)"
$ws.Range("E8").Formula = "=B8+1"

# ---------------------------------------------------------------------
# Row 9: continuation of the D8:D9 array formula; A9 is FALSE.
# ---------------------------------------------------------------------
$ws.Range("A9").Value = $false
$ws.Range("B9").Formula = "=C9"
$ws.Range("C9").Formula = "=D9"
$ws.Range("E9").Formula = "=B9+1"

# ---------------------------------------------------------------------
# Row 10: another circular loop, this time only B10 is a plain IF (not
# part of the array) while D10 starts a new array formula (D10:D11).
# ---------------------------------------------------------------------
$ws.Range("A10").Value = $true
$ws.Range("B10").Formula = "=IF(A10,C10,1)"
$ws.Range("C10").Formula = "=D10"
$ws.Range("D10:D11").FormulaArray = "=IF(A10:A11,B10:B11,1)"

# E10:E12 share one formula - E12 is cleared right after so it stays an
# empty (but present/formatted) cell, matching the trailing blank row.
$ws.Range("E10:E12").Formula = "=B10+1"
$ws.Range("E12").ClearContents()

# ---------------------------------------------------------------------
# Row 11: the loop resolves (A11 is FALSE) so every formula here settles
# to a concrete value instead of staying circular.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = $false
$ws.Range("B11").Formula = "=IF(A11,C11,1)"
$ws.Range("C11").Formula = "=D11"
$ws.Range("D11").Value = 1

# Move the saved cursor/selection off of the old D5 cell (the committed
# workbook no longer pins the view to a specific cell).
$ws.Range("A1").Select() | Out-Null

Write-Output "applied circular.xlsx row 8-12 edits"
